# Fixing geopoint in shared_table model
# Update schema.properties.* header labels on the "model" sheet so that
# latitude/longitude/altitude/accuracy point at their nested ".type" field.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

$ws.Range("E1").Value = "schema.properties.latitude.type"
$ws.Range("F1").Value = "schema.properties.longitude.type"
$ws.Range("G1").Value = "schema.properties.altitude.type"
$ws.Range("H1").Value = "schema.properties.accuracy.type"
